# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" column (E16:E22) with the new set of periods,
# now ordered from newest (2501) down to oldest (2407).
$ws.Range("E16").Value = "2501"
$ws.Range("E17").Value = "2412"
$ws.Range("E18").Value = "2411"
$ws.Range("E19").Value = "2410"
$ws.Range("E20").Value = "2409"
$ws.Range("E21").Value = "2408"
$ws.Range("E22").Value = "2407"

# The "Valor Mora" amount associated with period 2501 moves along with it to
# row 16, while row 22 (now period 2407) keeps the standard amount.
$ws.Range("F16").Value = 50266
$ws.Range("F22").Value = 52000
